# "Generate Report for Handoff"
# Updates the localization-status report: flips the status from
# "In Translation" to "Ready for handoff" and bumps the associated
# generation timestamps by one minute, then widens the Status-related
# columns so the new (longer) text isn't clipped.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-13 00:48:07"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-13 00:47:56"

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-13 00:48:07"

# --- Widen the Status columns to fit "Ready for handoff" -------------
# (target stored width ~17.216 chars; the host quantizes ColumnWidth to a
# 1/6-character pixel grid, so 16.3333... is the input that lands on the
# nearest reachable grid point, 17.1667 chars)
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333
